$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New draw-results row for 2025-12-05 (Pick 4), appended after the last
# existing row (79) -> becomes row 80.
$row = 80

# A80 ("2025-12-05") and C80 ("251205") look like a date / a plain number
# to Excel's type-inference, so force text storage via the "@" number
# format before writing the values, then clear the formatting again so
# the new cells don't end up with an explicit style (matches the look of
# every other row, which also carries no cell style).
$ws.Range("A80:E80").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2025-12-05"
$ws.Cells.Item($row, 2).Value = "Pick 4"
$ws.Cells.Item($row, 3).Value = "251205"
$ws.Cells.Item($row, 4).Value = "8-5-8-5"
$ws.Cells.Item($row, 5).Value = "2025-12-05T21:40:30.310+04:00"

$ws.Range("A80:E80").ClearFormats()
